# multi segment 3 changes
# --------------------------------------------------------------------------
# 1) "Booking Data" sheet: the "Return Date" column (C) is removed entirely
#    (Departure Date, previously D, slides left into C) and the (new) C2
#    departure-date value is bumped to 45637.
# --------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Booking Data")
$ws1.Range("C:C").EntireColumn.Delete()
$ws1.Range("C2").Value = 45637

# --------------------------------------------------------------------------
# 2) "multicity threesegment" sheet: the old "segment 1/2/3" block (columns
#    O..Z, 4 fields x 3 segments) is replaced by a simpler 3-leg layout
#    (Departure/Return/Date/Cabin Class per leg) re-using the columns
#    immediately after the shared trip fields, and the trailing unused
#    columns W:Z are dropped. A second (currently empty) data row is added
#    with just the Refundable flag.
# --------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("multicity threesegment")

# Drop the old segment-3 departure/return/date/class columns that are no
# longer used (W:Z) - this also shrinks the sheet dimension from Z to V.
$ws3.Range("W:Z").Delete()

# --- Row 1 (headers) --------------------------------------------------
# Columns A:C keep their existing meaning (Departure/Return/Departure Date).
# D used to be "Return Date" -> becomes "Cabin Class".
$ws3.Range("D1").Value = "Cabin Class"
# E:K keep the same relative meaning, just shifted from old F:L.
$ws3.Range("E1").Value = "Trip Type"
$ws3.Range("F1").Value = "Adults"
$ws3.Range("G1").Value = "Children"
$ws3.Range("H1").Value = "Infants"
$ws3.Range("I1").Value = "Direct Flight"
$ws3.Range("J1").Value = "Baggage"
$ws3.Range("K1").Value = "Key"
# L, M keep the same meaning as before (Language Code / Market Country Code).

# Leg 2 + leg 3 headers. The new unique shared strings must be created in
# this exact order (Q1 first, then N1/O1/P1, then R1/S1/T1/U1) to match the
# canonical shared-strings table ordering.
$ws3.Range("Q1").Value = "Cabin Class 2"
$ws3.Range("N1").Value = "Departure Location 2"
$ws3.Range("O1").Value = "Return Location 2"
$ws3.Range("P1").Value = "Departure Date 2"
$ws3.Range("R1").Value = "Departure Location 3"
$ws3.Range("S1").Value = "Return Location 3"
$ws3.Range("T1").Value = "Departure Date 3"
$ws3.Range("U1").Value = "Cabin Class 3"

# V1 "Refundable" header, re-using the bold/centered style already used for
# the equivalent header cell on the "Booking Data" sheet.
$ws1.Range("N1").Copy()
$ws3.Range("V1").PasteSpecial(-4122)
$ws3.Range("V1").Value = "Refundable"

# --- Row 2 (data) -------------------------------------------------------
$ws3.Range("D2").Value = "Y"
$ws3.Range("D2").NumberFormat = "General"

$ws3.Range("N2").Value = "DXB"
$ws3.Range("O2").Value = "CAI"
$ws3.Range("P2").Value = 45609
$ws3.Range("P2").NumberFormat = "d-mmm-yy"
$ws3.Range("Q2").Value = "Y"

$ws3.Range("R2").Value = "CAI"
$ws3.Range("S2").Value = "DXB"
$ws3.Range("T2").Value = 45611
$ws3.Range("T2").NumberFormat = "d-mmm-yy"
$ws3.Range("U2").Value = "Y"

$ws3.Range("V2").Value = $false

# --- Row 3 (new row, only the Refundable flag populated) ---------------
$ws3.Range("V3").Value = $false
